# The workbook originally holds a single sheet ("Tabelle1") with the
# timesheet content. This edit:
#   1. Inserts a new first sheet named "Logging" that stores small
#      config/carry-over values (used to remember row/column position
#      between sessions).
#   2. Renames the original sheet to "Timesheet" (it keeps all of its
#      existing content) and moves it to the second tab position.
#   3. Narrows the old multi-cell selection on the Timesheet sheet down
#      to a single cell (D12).
#   4. Leaves the Timesheet sheet active/selected (matching activeTab=1).

$wb = $excel.ActiveWorkbook

$originalSheet = $wb.Worksheets.Item(1)
$originalName = $originalSheet.Name

# Insert the new "Logging" sheet right before the existing sheet, so it
# becomes tab 1 and the original sheet becomes tab 2.
$logging = $wb.Worksheets.Add($originalSheet)
$logging.Name = "Logging"

# Re-fetch the original sheet by name (its index shifted to 2) and rename
# it to "Timesheet".
$timesheet = $wb.Worksheets.Item($originalName)
$timesheet.Name = "Timesheet"

# Populate the Logging sheet with the small carryover/row/column config
# table.
$logging.Range("B1").Value = "carryover"
$logging.Range("A2").Value = "row"
$logging.Range("B2").Value = 33
$logging.Range("A3").Value = "column"
$logging.Range("B3").Value = 10

# Narrow the Timesheet selection from D12:G31 down to just D12, and make
# the Timesheet sheet the active one.
$timesheet.Range("D12").Select() | Out-Null
$timesheet.Activate() | Out-Null
